$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 26
$ws.Range("B13").Value = 4642857.142857143
$ws.Range("B20").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B32").Value = 4642857.142857143
$ws.Range("B34").Value = 5342857.142857143
